$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("QuantitativeMetrics")

# Row 7: "Assertion validity" result changes from "no" to "yes", and the note is cleared.
$ws.Range("B7").Value = "yes"
$ws.Range("C7").ClearContents()
